# Fixed update to excel issue
#
# 1. Rename the generic "Requested quantity" header on the "Weekly Quantity"
#    and "Monthly Trend" sheets to more specific names.
# 2. Add a new "PO Forecast" sheet with forecast data
#    (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- 1. Rename the "Requested quantity" headers -----------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the "PO Forecast" sheet -----------------------------------------
# Duplicate the "Weekly Quantity" sheet (so it inherits the same sheetPr,
# pageMargins, header style and date-column style) and then overwrite its
# contents with the forecast data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"

# Extend the header formatting (bold font + border) from B1 across C1:D1.
$wsForecast.Range("B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

# Extend the date-column formatting down to row 43 (copied sheet only had 35
# rows of that style).
$wsForecast.Range("A35").Copy()
$wsForecast.Range("A36:A43").PasteSpecial(-4122)

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
    @(45312.99999999999, 508, 145.8098183648693, 892.7825178429489),
    @(45319.99999999999, 499, 108.5351727552197, 885.4973714335051),
    @(45326.99999999999, 489, 86.26438380504935, 886.8411011414785),
    @(45333.99999999999, 479, 96.78747476538096, 850.6203076796837),
    @(45347.99999999999, 459, 72.46969946524325, 839.2972522373504),
    @(45354.99999999999, 449, 92.30886614690128, 811.8210538415991),
    @(45361.99999999999, 439, 53.55109746185657, 812.7453622599968),
    @(45368.99999999999, 430, 48.281025509048, 818.6968236477281),
    @(45375.99999999999, 420, 19.50939532861754, 797.4604307820865),
    @(45382.99999999999, 410, -22.55947800121094, 801.5068749909087),
    @(45389.99999999999, 400, 33.45726121212674, 778.8054526408001),
    @(45403.99999999999, 380, 12.79225130361617, 755.8627254187193),
    @(45410.99999999999, 370, -43.92275861567186, 752.5309017946422),
    @(45424.99999999999, 351, -28.90206338481187, 746.9252543722971),
    @(45431.99999999999, 341, -17.05162522696716, 717.0589690107892),
    @(45438.99999999999, 331, -33.13055801650846, 712.8595845535627),
    @(45445.99999999999, 321, -52.00942880059735, 735.8314903610377),
    @(45459.99999999999, 301, -100.317534600808, 663.8396675127074),
    @(45466.99999999999, 291, -95.22387776150548, 676.671302277936),
    @(45473.99999999999, 282, -125.0599242051465, 653.0740043995708),
    @(45508.99999999999, 232, -155.2879231264015, 634.1637429501662),
    @(45522.99999999999, 212, -166.6597989179165, 610.0554944299174),
    @(45529.99999999999, 203, -205.4050412777973, 628.2492998968251),
    @(45536.99999999999, 193, -194.1484548597792, 585.2638765969939),
    @(45543.99999999999, 183, -188.8217694371561, 572.1630071077725),
    @(45550.99999999999, 173, -210.7186000038403, 544.7114384602068),
    @(45557.99999999999, 163, -221.3723763194249, 546.1274931102586),
    @(45571.99999999999, 143, -242.5867209130034, 535.2525543809338),
    @(45578.99999999999, 134, -243.2193633868621, 524.3498710247753),
    @(45585.99999999999, 124, -256.4454227334432, 509.5696280019298),
    @(45592.99999999999, 114, -276.5897221381379, 513.1809799671663),
    @(45599.99999999999, 104, -275.8229853388509, 467.8885914960937),
    @(45606.99999999999, 94, -313.0932545047821, 476.6300480915744),
    @(45620.99999999999, 74, -307.9767910293493, 472.4064108668519),
    @(45627.99999999999, 64, -310.8576390507456, 463.3051203147255),
    @(45634.99999999999, 55, -328.1111108190048, 410.3132637371115),
    @(45641.99999999999, 45, -333.0520965043895, 440.288913932659),
    @(45648.99999999999, 35, -353.4246054887493, 425.5501253092842),
    @(45655.99999999999, 25, -348.0745994576346, 446.8700911081787),
    @(45662.99999999999, 15, -361.6328951306547, 380.3307898342247),
    @(45669.99999999999, 5, -382.7468863562716, 399.8262093383179),
    @(45676.99999999999, 0, -409.404416594989, 382.4768442273719)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}
